$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '64.036.18'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '2.737.30'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue "D5" '569.56'
$ws.Range("E5").Value = '  -1.20%  '
Set-TextValue "D6" '159.36'
$ws.Range("E6").Value = '  +1.09%  '
Set-TextValue "D8" '0.597'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("E10").Value = '  +4.63%  '
Set-TextValue "D11" '5.72'
$ws.Range("E11").Value = '  -1.79%  '
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '3.221.69'
$ws.Range("E13").Value = '  -0.63%  '
Set-TextValue "D14" '26.72'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '63.618.82'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '2.746.39'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  -1.58%  '
Set-TextValue "D20" '354.02'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("E21").Value = '  -2.86%  '
$ws.Range("E22").Value = '  +0.08%  '
Set-TextValue "D23" '0.522'
$ws.Range("E23").Value = '  -5.37%  '
Set-TextValue "D24" '64.32'
$ws.Range("E24").Value = '  -2.81%  '
$ws.Range("E25").Value = '  +0.31%  '
Set-TextValue "D26" '1.00'
$ws.Range("E26").Value = '  +0.23%  '
Set-TextValue "D27" '8.43'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '0.0₃0909'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  +2.99%  '
$ws.Range("E31").Value = '  +7.90%  '
Set-TextValue "D32" '163.93'
$ws.Range("E32").Value = '  -3.09%  '
Set-TextValue "D33" '4.89'
$ws.Range("E33").Value = '  -0.61%  '
Set-TextValue "D34" '20.00'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("E35").Value = '  +1.59%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +0.97%  '
Set-TextValue "D38" '0.988'
$ws.Range("E38").Value = '  -0.74%  '
Set-TextValue "D39" '350.03'
$ws.Range("E39").Value = '  +5.91%  '
Set-TextValue "D40" '6.33'
$ws.Range("E40").Value = '  +2.69%  '
Set-TextValue "D41" '4.11'
$ws.Range("E41").Value = '  -1.24%  '
Set-TextValue "D42" '38.63'
$ws.Range("E42").Value = '  -1.63%  '
Set-TextValue "D43" '22.00'
$ws.Range("E43").Value = '  +1.39%  '
Set-TextValue "D44" '21.14'
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("E45").Value = '  -1.36%  '
Set-TextValue "D46" '0.625'
$ws.Range("E46").Value = '  -1.43%  '
Set-TextValue "D47" '134.58'
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("E48").Value = '  -1.06%  '
Set-TextValue "D49" '0.0249'
$ws.Range("E49").Value = '  -2.63%  '
Set-TextValue "D50" '0.999'
$ws.Range("E50").Value = '  -0.08%  '
Set-TextValue "D51" '11.06'
$ws.Range("E51").Value = '  +0.11%  '
